{"js": "// Remove the stale \"Ver no Jupiter Salvar em pdf Salvar em docx\" / \"\u00a9 2020 ...\"\n// footer paragraphs (plus the now-redundant blank paragraph that separated\n// them from the preceding page-break paragraph) from the end of the body.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\n// Locate the \"Ver no Jupiter...\" paragraph; the blank paragraph immediately\n// preceding it and the copyright paragraph immediately following it are\n// deleted together with it.\nconst items = paragraphs.items;\nlet jupiterIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === targets[0]) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex !== -1) {\n  const copyrightIndex = jupiterIndex + 1;\n  if (items[copyrightIndex] && items[copyrightIndex].text === targets[1]) {\n    items[copyrightIndex].delete();\n  }\n  items[jupiterIndex].delete();\n  const blankIndex = jupiterIndex - 1;\n  if (items[blankIndex] && items[blankIndex].text === \"\") {\n    items[blankIndex].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove the stale \"Ver no Jupiter Salvar em pdf Salvar em docx\" / \"(c) 2020 ...\"\n# footer paragraphs (plus the now-redundant blank paragraph that separated\n# them from the preceding page-break paragraph) from the end of the body.\n$d = $word.ActiveDocument\n\n$jupiterText = \"Ver no Jupiter Salvar em pdf Salvar em docx\"\n$copyrightText = [char]0x00A9 + \" 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n\n$jupiterIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $jupiterText) {\n        $jupiterIndex = $i\n        break\n    }\n}\n\nif ($jupiterIndex -ne -1) {\n    $copyrightIndex = $jupiterIndex + 1\n    $blankIndex = $jupiterIndex - 1\n\n    $copyrightMatches = $false\n    if ($copyrightIndex -le $d.Paragraphs.Count) {\n        $ct = $d.Paragraphs.Item($copyrightIndex).Range.Text.TrimEnd([char]13, [char]7)\n        if ($ct -eq $copyrightText) { $copyrightMatches = $true }\n    }\n\n    # Delete from the bottom up so earlier indices stay valid.\n    if ($copyrightMatches) {\n        $d.Paragraphs.Item($copyrightIndex).Range.Delete()\n    }\n    $d.Paragraphs.Item($jupiterIndex).Range.Delete()\n\n    if ($blankIndex -ge 1) {\n        $bt = $d.Paragraphs.Item($blankIndex).Range.Text.TrimEnd([char]13, [char]7)\n        if ($bt -eq \"\") {\n            $d.Paragraphs.Item($blankIndex).Range.Delete()\n        }\n    }\n}\n"}
